$p = $ppt.ActivePresentation
$cs = $p.SlideMaster.ColorScheme
for ($i=9; $i -le 12; $i++) {
  try {
    $c = $cs.Colors($i)
    Write-Output "$i : $($c.RGB)"
  } catch {
    Write-Output "$i : ERR $_"
  }
}
